$wb = $excel.ActiveWorkbook

$keySheet = $wb.Worksheets.Item("Key")
$dailySheet = $wb.Worksheets.Item("Daily Attendance Template")

# Clear the sample "day off" / ozeret pre-filled flags on the Key sheet,
# turning it back into a blank template (keep formatting/style where it existed).
$keySheet.Range("D2").ClearContents()
$keySheet.Range("E2").ClearContents()
$keySheet.Range("D3").ClearContents()
$keySheet.Range("E4").ClearContents()
$keySheet.Range("E5").ClearContents()
$keySheet.Range("D6").ClearContents()
$keySheet.Range("D7").ClearContents()

# Update the selection shown on the Key sheet.
$keySheet.Range("D2:F9").Select()

# Hide the Daily Attendance Template sheet and make Key the active/selected tab.
$dailySheet.Visible = $false
$keySheet.Activate()
